$d = $word.ActiveDocument

# Each <id>...</id> tag was previously split across three runs:
#   run1 (Courier New, color 7f6000): "<id>"
#   run2 (plain):                     "<oldvalue>"
#   run3 (Courier New, color 7f6000): "</id>"
# The edit collapses each triple into a single run containing the full
# "<id>newvalue</id>" text, formatted like the original surrounding tag
# markup (run1's formatting is inherited automatically by Word when the
# found Range's .Text is reassigned).

$replacements = @(
    @{ Old = "<id>p098v_a5</id>"; New = "<id>p098v_5</id>" },
    @{ Old = "<id>p099r_a2</id>"; New = "<id>p099r_2</id>" },
    @{ Old = "<id>p099r_a3</id>"; New = "<id>p099r_3</id>" }
)

foreach ($rep in $replacements) {
    $r = $d.Content
    $r.Find.ClearFormatting()
    $found = $r.Find.Execute($rep.Old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $r.Find.Found) {
        throw "Could not find text: $($rep.Old)"
    }
    $r.Text = $rep.New
}
